$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (202013211, Andrei Mercado) is a "Pending" alumni record that needs to be
# removed from this Approved Alumni sheet entirely.
$ws.Rows.Item(8).Delete()

# After the above deletion, the former row 9 (202119099, Biya Sungit) shifts up to
# become row 8, and the former rows 10 and 11 (Maria Benz, Rudolf Reindeer) shift up
# to become rows 9 and 10. Remove those two remaining "Pending" records as well.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()
